# Generate Report for Handback
# Updates the timestamps/priority recorded for the bc236758-*.md handback
# row (and the shared "Latest HO Xliff Generate Date" that row drives on
# the Overview sheet) to reflect a fresh report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the bc236758 / f28abb34 rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-16 02:14:56"
$wsOverview.Range("G5").Value = "2016-08-16 02:14:56"

# --- zh-cn sheet: Priority + Correspond Handoff/Handback Datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-16 02:14:51"
$wsZhCn.Range("H5").Value = "2016-08-16 02:14:51"
$wsZhCn.Range("K4").Value = "2016-08-16 02:15:16"
$wsZhCn.Range("K5").Value = "2016-08-16 02:15:16"

# --- de-de sheet: Priority + Correspond Handoff/Handback Datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-16 02:14:56"
$wsDeDe.Range("H5").Value = "2016-08-16 02:14:56"
$wsDeDe.Range("K4").Value = "2016-08-16 02:15:23"
$wsDeDe.Range("K5").Value = "2016-08-16 02:15:23"
